# Esquema general.docx - add "Método heurístico" section to the outline.
#
# 1) The bullet "Conclusión" (under "Extracción de características") gets an
#    extra sentence appended to it.
# 2) A whole new block of outline bullets is appended right after it, still
#    before the trailing empty paragraph / sectPr at the end of the document.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the "Conclusión" paragraph robustly (search by text instead of a
# hard-coded paragraph index).
# ---------------------------------------------------------------------------
$conclusion = $null
foreach ($para in $d.Paragraphs) {
    $t = $para.Range.Text
    $t = $t.TrimEnd([char]13, [char]7)
    if ($t -eq "Conclusión") {
        $conclusion = $para
        break
    }
}

if ($conclusion -eq $null) {
    throw "Could not find the 'Conclusión' paragraph"
}

# ---------------------------------------------------------------------------
# Append " -> mencionar parámetros del detector probados y los finalmente
# utilizados" right after the existing "Conclusión" text (before the
# paragraph mark).
# ---------------------------------------------------------------------------
$cr = $conclusion.Range
[void]$cr.MoveEnd(1, -1)
$cr.Collapse(0)
$cr.InsertAfter(" -> mencionar parámetros del detector probados y los finalmente utilizados")

# ---------------------------------------------------------------------------
# Insert the new "Método heurístico" outline block right after the
# "Conclusión" paragraph, one paragraph at a time, each one inheriting the
# list style/numbering ("Prrafodelista" / numId 1) and then getting its own
# indentation level set explicitly.
# ---------------------------------------------------------------------------
$newParagraphs = @(
    @{ Level = 1; Text = "Método heurístico" },
    @{ Level = 2; Text = "Objetivo -> Alimentar extracción de características (si se coge dos veces el mismo sello es problema de la ext. de caract. y no se trata en este punto)" },
    @{ Level = 2; Text = "Filtros" },
    @{ Level = 2; Text = "Resultados" },
    @{ Level = 1; Text = "Algoritmo final (con alimentación de sellos manual)" },
    @{ Level = 2; Text = "Crear base de datos de sellos" },
    @{ Level = 2; Text = "Crear matriz de acumulación de evidencias" },
    @{ Level = 2; Text = "Crear matriz de convolución de dicha acumulación de evidencias" },
    @{ Level = 2; Text = "Encontrar el máximo en dicha convolución -> Todos los puntos en dicha celda se consideran que forman parte del sello -> media de coordenadas = centro del sello" },
    @{ Level = 2; Text = "El sello cuyo ratio (matches dentro del sello)/(matches totales) sea mayor es el sello encontrado y se clasifica como tal -> Incluir pruebas de ratios" },
    @{ Level = 2; Text = "Se utilizan las coordenadas del centro y las dimensiones almacenadas del sello para eliminarlo del documento y se extrae la imagen sin sello para futuras operaciones." },
    @{ Level = 1; Text = "Resultados y posibles mejoras" }
)

$anchor = $conclusion
foreach ($item in $newParagraphs) {
    $anchor.Range.InsertParagraphAfter()
    $anchor = $anchor.Next()
    $anchor.Range.Text = $item.Text
    $anchor.Range.ListFormat.ListLevelNumber = $item.Level
}
